# Integrating mapping and grouping variables into updated interface.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add a new row 5, carrying over the old "farmsize" shortName/Levels/Labels/
#    definition/level/level_lab values under a new varName/label
#    (livestock_area / Land Area for Livestock).
$ws.Cells.Item(5, 1).Value = "livestock_area"
$ws.Cells.Item(5, 2).Value = "Land Area for Livestock"
$ws.Cells.Item(5, 3).Value = $ws.Cells.Item(4, 3).Value2
$ws.Cells.Item(5, 4).Value = $ws.Cells.Item(4, 4).Value2
$ws.Cells.Item(5, 5).Value = $ws.Cells.Item(4, 5).Value2
$ws.Cells.Item(5, 6).Value = $ws.Cells.Item(4, 6).Value2
$ws.Cells.Item(5, 7).Value = $ws.Cells.Item(4, 7).Value2
$ws.Cells.Item(5, 8).Value = $ws.Cells.Item(4, 8).Value2

# 2) Replace the "rural" disaggregate row (row 3) with a new "covid_shock"
#    disaggregate.
$ws.Cells.Item(3, 1).Value = "covid_shock"
$ws.Cells.Item(3, 3).Value = "COVID Impact"
$ws.Cells.Item(3, 2).Value = "Household Impacted by COVID"
$ws.Cells.Item(3, 4).Value = "1,2"
$ws.Cells.Item(3, 5).Value = "Yes,No"
$ws.Cells.Item(3, 6).Value = "1=Yes, 2=No"
$ws.Cells.Item(3, 7).Value = "Household"
$ws.Cells.Item(3, 8).Value = "Disaggregates"

# 3) Replace the "farmsize" disaggregate row (row 4) with a new "ag_comm"
#    disaggregate, reusing the 1,2 / Yes,No / 1=Yes, 2=No strings just
#    introduced above; clear the trailing level_lab cell (H4) entirely.
$ws.Cells.Item(4, 1).Value = "ag_comm"
$ws.Cells.Item(4, 2).Value = "Membership in an Agricultural Community"
$ws.Cells.Item(4, 3).Value = "Association Membership"
$ws.Cells.Item(4, 4).Value = "1,2"
$ws.Cells.Item(4, 5).Value = "Yes,No"
$ws.Cells.Item(4, 6).Value = "1=Yes, 2=No"
$ws.Cells.Item(4, 7).Value = "Household"
$ws.Cells.Item(4, 8).ClearContents()

# Match the saved selection/active-cell state seen in the edited workbook.
$ws.Range("E14").Select()
